$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 20

$ws.Range("A" + ($row - 1)).Copy()
$ws.Range("A" + $row).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = 42624.619722222225

$ws.Cells.Item($row, 2).Value = 32
$ws.Cells.Item($row, 3).Value = 49
$ws.Cells.Item($row, 4).Value = 50
$ws.Cells.Item($row, 5).Value = 49
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 9081
$ws.Cells.Item($row, 8).Value = 6589
$ws.Cells.Item($row, 9).Value = 437
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 61
$ws.Cells.Item($row, 12).Value = 16
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = "Named"

$wb.Save()
